$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
try {
  $sh = $s.Shapes.Item("Pentagon 197")
  Write-Host "Found Pentagon 197 directly: [$($sh.Name)] id=$($sh.Id)"
} catch {
  Write-Host "ERROR: $_"
}
